## email db sync done with validation
# Sync a handful of student records (email typo fix, a corrected surname,
# and a bounced/replaced email address) and re-validate the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- corrected / re-validated email for row 2 ---------------------------
$ws.Range("A2").Value = "alanwalker23@example.net"
# touching Font forces the cell through the validated-style path, same as
# the other re-validated records below
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Font.Color = 1

# --- corrected surname for row 7 -----------------------------------------
$ws.Range("B7").Value = "Justin Creamroll"

# --- corrected / re-validated email for row 8 ----------------------------
$ws.Range("A8").Value = "kirky@example.com"
$ws.Range("A8").Font.Name = "Calibri"
$ws.Range("A8").Font.Size = 11
$ws.Range("A8").Font.Color = 1

# --- widen the email column now that longer addresses are present --------
$ws.Columns.Item(1).ColumnWidth = 29.78

# --- tighten row heights across the header row and the validated block ---
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Range("A7:C50").RowHeight = 13.8

# --- leave the cursor/selection where the validation pass finished -------
$ws.Range("A7:C50").Select()
